# Updated Prediction values in column C for the PC Sun model retraining (Kahraman XGB results)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(4, 3).Value = 0.246
$ws.Cells.Item(5, 3).Value = 0.876
$ws.Cells.Item(6, 3).Value = 1.474
$ws.Cells.Item(7, 3).Value = 1.538
$ws.Cells.Item(9, 3).Value = 1.427
$ws.Cells.Item(10, 3).Value = 1.088
$ws.Cells.Item(11, 3).Value = 0.638
$ws.Cells.Item(27, 3).Value = 0.024
$ws.Cells.Item(28, 3).Value = 0.103
$ws.Cells.Item(29, 3).Value = 0.216
$ws.Cells.Item(30, 3).Value = 0.353
$ws.Cells.Item(31, 3).Value = 0.401
$ws.Cells.Item(32, 3).Value = 0.452
$ws.Cells.Item(33, 3).Value = 0.385
$ws.Cells.Item(34, 3).Value = 0.309
$ws.Cells.Item(35, 3).Value = 0.15
$ws.Cells.Item(36, 3).Value = 0.052
$ws.Cells.Item(52, 3).Value = 0.434
$ws.Cells.Item(53, 3).Value = 1.145
$ws.Cells.Item(54, 3).Value = 1.657
$ws.Cells.Item(55, 3).Value = 2.024
$ws.Cells.Item(56, 3).Value = 2.153
$ws.Cells.Item(57, 3).Value = 1.945
$ws.Cells.Item(58, 3).Value = 1.479
$ws.Cells.Item(59, 3).Value = 0.754
$ws.Cells.Item(65, 3).Value = 0.01
$ws.Cells.Item(71, 3).Value = 0.011
$ws.Cells.Item(72, 3).Value = 0.013
$ws.Cells.Item(73, 3).Value = 0.012
$ws.Cells.Item(74, 3).Value = 0.012
$ws.Cells.Item(75, 3).Value = 0.058
$ws.Cells.Item(76, 3).Value = 0.607
$ws.Cells.Item(77, 3).Value = 1.639
$ws.Cells.Item(78, 3).Value = 2.127
$ws.Cells.Item(79, 3).Value = 2.904
$ws.Cells.Item(80, 3).Value = 2.956
$ws.Cells.Item(81, 3).Value = 2.255
$ws.Cells.Item(82, 3).Value = 1.789
$ws.Cells.Item(95, 3).Value = 0.011
$ws.Cells.Item(96, 3).Value = 0.014
$ws.Cells.Item(97, 3).Value = 0.014
$ws.Cells.Item(99, 3).Value = 0.054
$ws.Cells.Item(100, 3).Value = 0.5610000000000001
$ws.Cells.Item(101, 3).Value = 1.589
$ws.Cells.Item(102, 3).Value = 2.051
$ws.Cells.Item(103, 3).Value = 2.289
$ws.Cells.Item(104, 3).Value = 2.319
$ws.Cells.Item(105, 3).Value = 2.088
$ws.Cells.Item(106, 3).Value = 1.613
$ws.Cells.Item(107, 3).Value = 0.679
$ws.Cells.Item(108, 3).Value = 0.139
$ws.Cells.Item(109, 3).Value = 0.012
$ws.Cells.Item(112, 3).Value = 0.02
$ws.Cells.Item(166, 3).Value = 0.012
$ws.Cells.Item(170, 3).Value = 0.018